# Add support for 'CELL("contents",[reference])' command.
# Mirrors the author's manual edit: for each informational column
# (B..K) on row 2, enter the array formula
#   =CELL("contents", <col>2)
# into the corresponding cell of the new row 22, array-entered
# (Ctrl+Shift+Enter) exactly like the existing A5 na() array formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

foreach ($col in $cols) {
    $targetCell = $col + "22"
    $sourceCell = $col + "2"
    $ws.Range($targetCell).FormulaArray = '=CELL("contents", ' + $sourceCell + ')'
}

# Leave the selection where the author's cursor ended up after
# array-entering the last formula (one cell to the right of the block).
$ws.Range("L22").Select() | Out-Null
